# Update column F (dSF) values for specific rows in Sheet1
# per the commit "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = 2
$ws.Range("F10").Value = -4
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 4
$ws.Range("F14").Value = 2
$ws.Range("F15").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = -1
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = -4
$ws.Range("F21").Value = -1
$ws.Range("F22").Value = -2
$ws.Range("F26").Value = 4
$ws.Range("F27").Value = 1
$ws.Range("F28").Value = 1
